# Update the Big Data exercise row (row 18) on the class_schedule sheet:
#  - Append a third bullet to the "Do Before Class" cell (E18)
#  - Add the new in-class exercise link to F18
#  - Update the sheet view / selection to reflect the edited cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E18").Value = '- `What is Big Data? <what_is_big_data.ipynb>`_' + [char]10 + '- `Strategies for Big Data <big_data_strategies.ipynb>`_' + [char]10 + '- Download the dataset linked at the top of the linked exercise.'
$ws.Range("F18").Value = '`Link <exercises/Exercise_bigdata.ipynb>`_'

$ws.Rows.Item(18).RowHeight = 68

$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("F18").Select()
